$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "1.003") must be forced
# to Text first, otherwise Excel auto-converts the string to a floating-point
# number (losing the exact text representation used throughout this sheet).
$numericLookingCells = @("D4","D5","D7","D8","D9","D10","D11","D12","D13","D14","D15","D18","D19","D20","D21","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "23.876.50"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "1.654.52"
$ws.Range("E3").Value = "  +2.50%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "308.76"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.3888"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "0.3838"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("D9").Value = "50.93"
$ws.Range("E9").Value = "  +5.05%  "
$ws.Range("D10").Value = "1.351"
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "0.08462"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").Value = "23.88"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").Value = "7.141"
$ws.Range("E14").Value = "  +2.47%  "
$ws.Range("D15").Value = "7.876"
$ws.Range("E15").Value = "  +6.55%  "
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("D17").Value = "1.651.30"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").Value = "94.69"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("D19").Value = "0.07007"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").Value = "19.76"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").Value = "6.902"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "13.58"
$ws.Range("E23").Value = "  +1.82%  "
$ws.Range("D24").Value = "23.869.26"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").Value = "2.489"
$ws.Range("E25").Value = "  +3.07%  "
$ws.Range("D26").Value = "3.043"
$ws.Range("E26").Value = "  +8.19%  "
$ws.Range("D27").Value = "22.05"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").Value = "152.46"
$ws.Range("E28").Value = "  -2.77%  "
$ws.Range("D29").Value = "5.414"
$ws.Range("E29").Value = "  +3.48%  "
$ws.Range("D30").Value = "139.11"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").Value = "7.754"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").Value = "2.496"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Value = "1.833.67"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("D34").Value = "1.028"
$ws.Range("E34").Value = "  +8.60%  "
$ws.Range("D35").Value = "0.08047"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").Value = "0.02956"
$ws.Range("E36").Value = "  +4.07%  "
$ws.Range("D37").Value = "10.96"
$ws.Range("E37").Value = "  +5.27%  "
$ws.Range("D38").Value = "6.650"
$ws.Range("E38").Value = "  +2.13%  "
$ws.Range("D39").Value = "0.2679"
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("D40").Value = "0.09108"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.7524"
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "13.45"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("D43").Value = "1.420"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "16.25"
$ws.Range("E44").Value = "  +3.34%  "
$ws.Range("D45").Value = "0.6942"
$ws.Range("E45").Value = "  +2.37%  "
$ws.Range("D46").Value = "2.459"
$ws.Range("E46").Value = "  +1.54%  "
$ws.Range("D47").Value = "4.072"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "0.08254"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").Value = "134.35"
$ws.Range("E50").Value = "  +2.00%  "
$ws.Range("D51").Value = "1.232"
$ws.Range("E51").Value = "  +8.45%  "

foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
